# Update cryptocurrency price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "43.909.81"
$ws.Cells.Item(2, 5).Value = "  -0.93%  "
$ws.Cells.Item(3, 4).Value = "2.356.62"
$ws.Cells.Item(3, 5).Value = "  -0.32%  "
$ws.Cells.Item(4, 5).Value = "  -0.03%  "
$ws.Cells.Item(5, 4).Value = "'0.679"
$ws.Cells.Item(5, 5).Value = "  +0.51%  "
$ws.Cells.Item(6, 4).Value = "'239.31"
$ws.Cells.Item(6, 5).Value = "  +0.11%  "
$ws.Cells.Item(7, 4).Value = "'73.68"
$ws.Cells.Item(7, 5).Value = "  -0.39%  "
$ws.Cells.Item(8, 5).Value = "  -0.04%  "
$ws.Cells.Item(9, 4).Value = "'0.590"
$ws.Cells.Item(9, 5).Value = "  +6.62%  "
$ws.Cells.Item(10, 5).Value = "  -1.46%  "
$ws.Cells.Item(11, 4).Value = "'57.27"
$ws.Cells.Item(11, 5).Value = "  -0.26%  "
$ws.Cells.Item(12, 4).Value = "'32.32"
$ws.Cells.Item(12, 5).Value = "  +6.38%  "
$ws.Cells.Item(13, 5).Value = "  +0.63%  "
$ws.Cells.Item(14, 4).Value = "'7.21"
$ws.Cells.Item(14, 5).Value = "  +4.69%  "
$ws.Cells.Item(15, 4).Value = "2.707.18"
$ws.Cells.Item(15, 5).Value = "  -0.17%  "
$ws.Cells.Item(16, 4).Value = "'16.46"
$ws.Cells.Item(16, 5).Value = "  -2.54%  "
$ws.Cells.Item(17, 4).Value = "'0.898"
$ws.Cells.Item(17, 5).Value = "  -0.98%  "
$ws.Cells.Item(18, 4).Value = "2.363.16"
$ws.Cells.Item(18, 5).Value = "  +0.03%  "
$ws.Cells.Item(19, 4).Value = "43.790.12"
$ws.Cells.Item(19, 5).Value = "  -1.50%  "
$ws.Cells.Item(20, 4).Value = "'6.89"
$ws.Cells.Item(20, 5).Value = "  +6.04%  "
$ws.Cells.Item(21, 5).Value = "  -1.20%  "
$ws.Cells.Item(22, 4).Value = "'77.41"
$ws.Cells.Item(22, 5).Value = "  -0.06%  "
$ws.Cells.Item(23, 4).Value = "'257.04"
$ws.Cells.Item(23, 5).Value = "  +0.90%  "
$ws.Cells.Item(24, 5).Value = "  +23.01%  "
$ws.Cells.Item(25, 5).Value = "  -0.06%  "
$ws.Cells.Item(26, 5).Value = "  -5.16%  "
$ws.Cells.Item(27, 5).Value = "  -1.80%  "
$ws.Cells.Item(28, 4).Value = "'10.77"
$ws.Cells.Item(28, 5).Value = "  +3.78%  "
$ws.Cells.Item(29, 5).Value = "  +1.52%  "
$ws.Cells.Item(30, 4).Value = "'22.81"
$ws.Cells.Item(30, 5).Value = "  +0.65%  "
$ws.Cells.Item(31, 4).Value = "'175.71"
$ws.Cells.Item(31, 5).Value = "  +1.06%  "
$ws.Cells.Item(32, 5).Value = "  -1.43%  "
$ws.Cells.Item(33, 5).Value = "  +2.17%  "
$ws.Cells.Item(34, 4).Value = "'0.0755"
$ws.Cells.Item(34, 5).Value = "  +1.54%  "
$ws.Cells.Item(35, 4).Value = "'5.58"
$ws.Cells.Item(35, 5).Value = "  +7.07%  "
$ws.Cells.Item(36, 4).Value = "'5.17"
$ws.Cells.Item(36, 5).Value = "  -0.74%  "
$ws.Cells.Item(37, 5).Value = "  -3.81%  "
$ws.Cells.Item(38, 5).Value = "  -2.43%  "
$ws.Cells.Item(39, 5).Value = "  -3.69%  "
$ws.Cells.Item(40, 5).Value = "  +2.29%  "
$ws.Cells.Item(41, 5).Value = "  +11.62%  "
$ws.Cells.Item(42, 4).Value = "'0.202"
$ws.Cells.Item(42, 5).Value = "  +8.99%  "
$ws.Cells.Item(43, 5).Value = "  +1.50%  "
$ws.Cells.Item(44, 4).Value = "'18.81"
$ws.Cells.Item(44, 5).Value = "  -2.75%  "
$ws.Cells.Item(45, 5).Value = "  -0.06%  "
$ws.Cells.Item(46, 4).Value = "'59.61"
$ws.Cells.Item(46, 5).Value = "  +13.28%  "
$ws.Cells.Item(47, 4).Value = "'4.76"
$ws.Cells.Item(47, 5).Value = "  +6.60%  "
$ws.Cells.Item(48, 5).Value = "  +3.62%  "
$ws.Cells.Item(49, 5).Value = "  -0.90%  "
$ws.Cells.Item(50, 4).Value = "'100.53"
$ws.Cells.Item(50, 5).Value = "  +1.31%  "
$ws.Cells.Item(51, 5).Value = "  -0.71%  "
